$d = $word.ActiveDocument

$replacements = @(
    @("35×84=2940", "98×53=5194"),
    @("49×13=637", "28×89=2492"),
    @("96×12=1152", "11×21=231"),
    @("51×49=2499", "67×38=2546"),
    @("66×81=5346", "35×86=3010"),
    @("43×97=4171", "15×89=1335"),
    @("55×36=1980", "96×28=2688"),
    @("18×86=1548", "68×25=1700"),
    @("42×24=1008", "19×85=1615"),
    @("83×39=3237", "20×26=520"),
    @("76×80=6080", "29×36=1044"),
    @("95×83=7885", "21×35=735"),
    @("93×40=3720", "86×78=6708"),
    @("83×34=2822", "14×97=1358"),
    @("40×19=760", "98×55=5390"),
    @("87×92=8004", "74×29=2146"),
    @("81×97=7857", "59×52=3068"),
    @("43×82=3526", "59×13=767"),
    @("50×14=700", "64×18=1152"),
    @("31×46=1426", "53×44=2332"),
    @("17×44=748", "27×83=2241"),
    @("36×29=1044", "28×57=1596"),
    @("13×98=1274", "23×25=575"),
    @("83×41=3403", "44×78=3432"),
    @("88×32=2816", "45×98=4410")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
